$wb = $excel.ActiveWorkbook

# --- 1. Update the Date value on the Metadata sheet ---
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B8").Value = "2023-03-22T16:32:25+01:00"

# --- 2. Clear cells that held the removed "1MB" note on the Elements sheet ---
$wsElem = $wb.Worksheets.Item("Elements")
$wsElem.Range("N15").Value = ""
$wsElem.Range("N16").Value = ""

# --- 3. Clear the "ele-1" Condition(s) values in column AI for the listed rows ---
$rows = @(4,6,7,9,10,12,13,14,15,16,17,18)
foreach ($r in $rows) {
    $wsElem.Range("AI$r").Value = ""
}

# --- 4. Fix casing of "N/A" -> "n/a" in column AK for rows 9 and 12 ---
$wsElem.Range("AK9").Value = "n/a"
$wsElem.Range("AK12").Value = "n/a"
